# Populate the "November" worksheet's borrower/lender statistics that were
# previously left blank (B:G, rows 2-55). Values below are the literal
# per-library totals for November 2023 (Net lender / Net borrower /
# difference), plus the derived "borrowed more / lent more" note and the
# borrowed:lent ratio text - matching the pattern already used on every
# other populated month sheet (January-October, December stays blank).
#
# Once these literals are written, the "Yearly totals" sheet recalculates
# automatically (its B/C/D/G cells are formulas that sum every month sheet,
# including November), so it needs no direct edits here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("November")

$rows = @(
    @{ Row=2; B=1468; C=1211; D=257; E=1; F=0 },
    @{ Row=3; B=450; C=431; D=19; E=1; F=0 },
    @{ Row=4; B=1092; C=1418; D=-326; E=0; F=1 },
    @{ Row=5; B=22; C=155; D=-133; E=0; F=1 },
    @{ Row=6; B=1141; C=1290; D=-149; E=0; F=1 },
    @{ Row=7; B=243; C=151; D=92; E=1; F=0 },
    @{ Row=8; B=70; C=179; D=-109; E=0; F=1 },
    @{ Row=9; B=70; C=72; D=-2; E=0; F=1 },
    @{ Row=10; B=2; C=54; D=-52; E=0; F=1 },
    @{ Row=11; B=0; C=0; D=0; E=0; F=0 },
    @{ Row=12; B=41; C=25; D=16; E=1; F=0 },
    @{ Row=13; B=182; C=83; D=99; E=1; F=0 },
    @{ Row=14; B=128; C=262; D=-134; E=0; F=1 },
    @{ Row=15; B=97; C=90; D=7; E=1; F=0 },
    @{ Row=16; B=19; C=141; D=-122; E=0; F=1 },
    @{ Row=17; B=593; C=452; D=141; E=1; F=0 },
    @{ Row=18; B=67; C=109; D=-42; E=0; F=1 },
    @{ Row=19; B=629; C=379; D=250; E=1; F=0 },
    @{ Row=20; B=4; C=68; D=-64; E=0; F=1 },
    @{ Row=21; B=405; C=422; D=-17; E=0; F=1 },
    @{ Row=22; B=31; C=128; D=-97; E=0; F=1 },
    @{ Row=23; B=766; C=454; D=312; E=1; F=0 },
    @{ Row=24; B=2023; C=1209; D=814; E=1; F=0 },
    @{ Row=25; B=208; C=368; D=-160; E=0; F=1 },
    @{ Row=26; B=0; C=0; D=0; E=0; F=0 },
    @{ Row=27; B=278; C=209; D=69; E=1; F=0 },
    @{ Row=28; B=44; C=99; D=-55; E=0; F=1 },
    @{ Row=29; B=570; C=488; D=82; E=1; F=0 },
    @{ Row=30; B=24; C=32; D=-8; E=0; F=1 },
    @{ Row=31; B=86; C=293; D=-207; E=0; F=1 },
    @{ Row=32; B=419; C=558; D=-139; E=0; F=1 },
    @{ Row=33; B=344; C=417; D=-73; E=0; F=1 },
    @{ Row=34; B=216; C=129; D=87; E=1; F=0 },
    @{ Row=35; B=853; C=958; D=-105; E=0; F=1 },
    @{ Row=36; B=192; C=420; D=-228; E=0; F=1 },
    @{ Row=37; B=435; C=337; D=98; E=1; F=0 },
    @{ Row=38; B=19; C=175; D=-156; E=0; F=1 },
    @{ Row=39; B=35; C=76; D=-41; E=0; F=1 },
    @{ Row=40; B=86; C=118; D=-32; E=0; F=1 },
    @{ Row=41; B=2; C=43; D=-41; E=0; F=1 },
    @{ Row=42; B=15; C=21; D=-6; E=0; F=1 },
    @{ Row=43; B=0; C=0; D=0; E=0; F=0 },
    @{ Row=44; B=74; C=89; D=-15; E=0; F=1 },
    @{ Row=45; B=77; C=173; D=-96; E=0; F=1 },
    @{ Row=46; B=507; C=524; D=-17; E=0; F=1 },
    @{ Row=47; B=930; C=724; D=206; E=1; F=0 },
    @{ Row=48; B=229; C=572; D=-343; E=0; F=1 },
    @{ Row=49; B=457; C=236; D=221; E=1; F=0 },
    @{ Row=50; B=1043; C=505; D=538; E=1; F=0 },
    @{ Row=51; B=163; C=186; D=-23; E=0; F=1 },
    @{ Row=52; B=458; C=475; D=-17; E=0; F=1 },
    @{ Row=53; B=100; C=247; D=-147; E=0; F=1 },
    @{ Row=54; B=3; C=255; D=-252; E=0; F=1 },
    @{ Row=55; B=340; C=240; D=100; E=1; F=0 }
)

$noteBorrowed = "We borrowerd more than we lent"
$noteLent     = "We lent more than we borrowed"

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 2).Value = $r.B   # B: Other Next library materials checked out at our library
    $ws.Cells.Item($row, 3).Value = $r.C   # C: Our materials checked out at other Next libraries
    $ws.Cells.Item($row, 4).Value = $r.D   # D: Net (B - C)

    if ($r.E -eq 1) {
        $ws.Cells.Item($row, 5).Value = $noteBorrowed
    }
    if ($r.F -eq 1) {
        $ws.Cells.Item($row, 6).Value = $noteLent
    }

    if ($r.C -ne 0) {
        $ratio = [Math]::Round($r.B / $r.C, 2)
        $ws.Cells.Item($row, 7).Value = "$ratio : 1"
    }
}
